# lines_states.xlsx update: add line7/line8 entries, renumber subsequent
# "extr" rows, add two new rows (16 and 17) and refresh several
# from_bus/to_bus/in_service values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update labels in column B for rows 8-15 (the "line"/"extr" series
#     shifts down by two positions once line7/line8 are inserted) -------
$ws.Range("B8").Value2  = "line7"
$ws.Range("B9").Value2  = "line8"
$ws.Range("B10").Value2 = "extr1"
$ws.Range("B11").Value2 = "extr2"
$ws.Range("B12").Value2 = "extr3"
$ws.Range("B13").Value2 = "extr4"
$ws.Range("B14").Value2 = "extr5"
$ws.Range("B15").Value2 = "extr6"

# --- Update from_bus (C) / to_bus (D) / in_service (E) values ----------
$ws.Range("C8").Value2  = 14
$ws.Range("D8").Value2  = 11

$ws.Range("C9").Value2  = 16
$ws.Range("E9").Value2  = $true

$ws.Range("C10").Value2 = 5
$ws.Range("D10").Value2 = 12

$ws.Range("C11").Value2 = 5
$ws.Range("D11").Value2 = 9
$ws.Range("E11").Value2 = $true

$ws.Range("C12").Value2 = 10
$ws.Range("E12").Value2 = $true

$ws.Range("D13").Value2 = 8

$ws.Range("C14").Value2 = 9
$ws.Range("D14").Value2 = 11

$ws.Range("C15").Value2 = 7
$ws.Range("D15").Value2 = 11
$ws.Range("E15").Value2 = $true

# --- Add two new rows (16 and 17) --------------------------------------
$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "extr7"
$ws.Range("C16").Value2 = 5
$ws.Range("D16").Value2 = 7
$ws.Range("E16").Value2 = $true

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "extr8"
$ws.Range("C17").Value2 = 8
$ws.Range("D17").Value2 = 5
$ws.Range("E17").Value2 = $true

# Column A on the new rows should carry the same formatting (bold,
# centered, thin box border) used by the rest of the A column.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
